# Update "Datos actualizados" timestamp string (row 1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 9 de Octubre de 2020 a las 10:34"

# Refresh COVID-19 numbers for Rusia (row 7)
$ws.Range("B7").Value = 1272238
$ws.Range("C7").Value = 12126
$ws.Range("D7").Value = 1009421
$ws.Range("E7").Value = 240560
$ws.Range("G7").Value = 201
$ws.Range("H7").Value = 22257

# Refresh COVID-19 numbers for Filipinas (row 22)
$ws.Range("B22").Value = 334770
$ws.Range("C22").Value = 2996
$ws.Range("D22").Value = 275307
$ws.Range("E22").Value = 53311
$ws.Range("G22").Value = 83
$ws.Range("H22").Value = 6152

# Refresh COVID-19 numbers for Indonesia (row 24)
$ws.Range("B24").Value = 324658
$ws.Range("C24").Value = 4094
$ws.Range("D24").Value = 247667
$ws.Range("E24").Value = 65314
$ws.Range("G24").Value = 97
$ws.Range("H24").Value = 11677

# Refresh COVID-19 numbers for Singapur (row 62)
$ws.Range("B62").Value = 57859
$ws.Range("C62").Value = 10
$ws.Range("E62").Value = 164

# Lituania overtakes Suazilandia in the ranking (rows 123-124)
$ws.Range("A123").Value = "Lituania"
$ws.Range("B123").Value = 5758
$ws.Range("C123").Value = 133
$ws.Range("D123").Value = 2722
$ws.Range("E123").Value = 2933
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = 103

$ws.Range("A124").Value = "Suazilandia"
$ws.Range("B124").Value = 5632
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 5231
$ws.Range("E124").Value = 288
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 113

# Refresh COVID-19 numbers for Estonia (row 140)
$ws.Range("B140").Value = 3809
$ws.Range("C140").Value = 49
$ws.Range("D140").Value = 2906
$ws.Range("E140").Value = 836

# Letonia climbs above Polinesia Francesa, Benin, Guinea-Bisau and Belice (rows 151-155)
$ws.Range("A151").Value = "Letonia"
$ws.Range("B151").Value = 2507
$ws.Range("C151").Value = 137
$ws.Range("D151").Value = 1322
$ws.Range("E151").Value = 1145
$ws.Range("H151").Value = 40

$ws.Range("A152").Value = "Polinesia Francesa"
$ws.Range("B152").Value = 2420
$ws.Range("D152").Value = 1857
$ws.Range("E152").Value = 553
$ws.Range("H152").Value = 10

$ws.Range("A153").Value = "Benin"
$ws.Range("B153").Value = 2411
$ws.Range("D153").Value = 1973
$ws.Range("E153").Value = 397
$ws.Range("H153").Value = 41

$ws.Range("A154").Value = "Guinea-Bisau"
$ws.Range("B154").Value = 2385
$ws.Range("C154").Value = 0
$ws.Range("D154").Value = 1728
$ws.Range("E154").Value = 617
$ws.Range("H154").Value = 40

$ws.Range("A155").Value = "Belice"
$ws.Range("B155").Value = 2373
$ws.Range("C155").Value = 63
$ws.Range("D155").Value = 1459
$ws.Range("E155").Value = 880
$ws.Range("H155").Value = 34
